# Applies the "READ ME.docx" revision:
#   * stamps the inline picture's <wp:inline> with wp14:anchorId / wp14:editId
#     (the id pair Word mints whenever it rewrites a drawing's XML)
#   * rewrites step "2." to reference the new SRM_Package data folder
#   * rewrites step "3." to describe editing srm.markdown.R
#   * rewrites step "4." into the fuller walkthrough of editing line 30 and
#     appends new step "5." describing where the report is saved
#
# The new wording introduces w:proofErr (spell-check) markers and a few run
# splits that the Range/Find object model cannot author directly, so each
# paragraph is replaced wholesale via Range.InsertXML with the exact target
# WordprocessingML for that paragraph (this also lets the picture paragraph
# gain its two new attributes without disturbing anything else inside it).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Picture paragraph -- add wp14:anchorId/wp14:editId to <wp:inline>.
# ---------------------------------------------------------------------
$picturePara = $d.Paragraphs.Item(7)
if ($picturePara.Range.InlineShapes.Count -ne 1) {
    throw "Expected the picture paragraph at index 7"
}
$picturePayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mo="http://schemas.microsoft.com/office/mac/office/2008/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:mv="urn:schemas-microsoft-com:mac:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p w:rsidR="00C93138" w:rsidRDefault="00C93138"><w:r><w:rPr><w:noProof/><w:lang w:val="en-US" w:eastAsia="en-US"/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="388343DD" wp14:editId="1ECF4B74"><wp:extent cx="5270500" cy="1689735"/><wp:effectExtent l="0" t="0" r="12700" b="12065"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="SRM_Op.png"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5270500" cy="1689735"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$picturePara.Range.InsertXML($picturePayload)

# ---------------------------------------------------------------------
# 2. "2.  Save Excel data to the SRM Review folder as..." paragraph.
# ---------------------------------------------------------------------
$step2Para = $d.Paragraphs.Item(9)
if ($step2Para.Range.Text -notmatch "^2\.") {
    throw "Expected step 2 paragraph at index 9, got: $($step2Para.Range.Text)"
}
$step2Payload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mo="http://schemas.microsoft.com/office/mac/office/2008/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:mv="urn:schemas-microsoft-com:mac:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p w:rsidR="00C93138" w:rsidRDefault="00C93138"><w:r><w:t xml:space="preserve">2.  Save Excel data to the </w:t></w:r><w:r><w:t xml:space="preserve">data folder in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SRM_Package</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> set</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$step2Para.Range.InsertXML($step2Payload)

# ---------------------------------------------------------------------
# 3. "3.  Overwrite the current file, if asked." paragraph.
# ---------------------------------------------------------------------
$step3Para = $d.Paragraphs.Item(11)
if ($step3Para.Range.Text -notmatch "^3\.") {
    throw "Expected step 3 paragraph at index 11, got: $($step3Para.Range.Text)"
}
$step3Payload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mo="http://schemas.microsoft.com/office/mac/office/2008/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:mv="urn:schemas-microsoft-com:mac:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p w:rsidR="00C93138" w:rsidRDefault="00C93138"><w:r><w:t xml:space="preserve">3.  </w:t></w:r><w:r><w:t>Open “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>srm.markdown.R</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” and change filename to import data file</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$step3Para.Range.InsertXML($step3Payload)

# ---------------------------------------------------------------------
# 4. "4. Run Control Chart Review.Rmd" paragraph -> six new paragraphs
#    (step 4, the indented code line, the explanation, a blank line, and
#    the new step 5). Stop one character short of the paragraph's end so
#    the document's final paragraph mark is reused instead of InsertXML
#    appending a spurious blank paragraph before the sectPr.
# ---------------------------------------------------------------------
$step4Para = $d.Paragraphs.Item(13)
if ($step4Para.Range.Text -notmatch "^4\.") {
    throw "Expected step 4 paragraph at index 13, got: $($step4Para.Range.Text)"
}
$step4Range = $d.Range($step4Para.Range.Start, $step4Para.Range.End - 1)
$step4Payload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mo="http://schemas.microsoft.com/office/mac/office/2008/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:mv="urn:schemas-microsoft-com:mac:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p w:rsidR="00C93138" w:rsidRDefault="00C93138"><w:r><w:t xml:space="preserve">4. </w:t></w:r><w:r><w:t xml:space="preserve">Change the line 30  </w:t></w:r></w:p><w:p w:rsidR="00C93138" w:rsidRDefault="00C93138"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>srm_report</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>data.raw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>max.pts</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 200, points = 30, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>doc_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>docx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”)</w:t></w:r><w:r><w:t>”</w:t></w:r></w:p><w:p w:rsidR="00C93138" w:rsidRDefault="00C93138"><w:r><w:t>to reflect the number of data points, how many will be used to set the control limits and what type of output is desired.  (“html”, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>docx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, or “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pdf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”.).</w:t></w:r></w:p><w:p w:rsidR="00C93138" w:rsidRDefault="00C93138"/><w:p w:rsidR="00C93138" w:rsidRDefault="00C93138"><w:r><w:t>5.  Resultant report will be saved to data/Output folder.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$step4Range.InsertXML($step4Payload)

Write-Host "READ ME.docx steps updated."
